$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 7 ("Incendio" / "Se cuentan con extintores...") is removed entirely;
# the rows below it shift up by one.
$ws.Rows.Item(7).Delete()

# Row 5 (Robo de equipos...): probability lowered 0.2 -> 0.1
$ws.Range("C5").Value = 0.1

# Row 6 (Inundacion...) now also affects "Todas" areas, and the mitigation text was reworded
$ws.Range("B6").Value = "Todas"
$ws.Range("D6").Value = "Ubicar los servidires a 1 m. de altura, a salvo de posibles inundaciones. Además equipar la sala con alcantarillas para desagotar rápidamente cualquier fuga de agua."

# Row 7 (was "Corte de energia electrica") now reworded, affects "Todas" areas, probability changed, and
# the cell should wrap its text like the row above it
$ws.Range("A7").Value = "Corte de energía eléctrica debido a fallas por parte del proveedor"
$ws.Range("A7").WrapText = $true
$ws.Range("B7").Value = "Todas"
$ws.Range("C7").Value = 0.3

# The remaining risk rows were reordered
$ws.Range("A11").Value = "Falla total o parcial del cableado"
$ws.Range("A12").Value = "Problemas en el cableado eléctrico de las Estaciones de trabajo"
$ws.Range("A15").Value = "Problemas con los recursos compartidos de la red"
$ws.Range("A16").Value = "Pérdida total de un servidor"

$ws.Range("C4").Select()
